# Add a new User Story row to the Compare-genie UserStories worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row is appended right after the last existing row (row 28 -> row 29).
$prevRow = 28
$newRow = 29

# Carry over the formatting of the previous (last) row's populated cells
# (A28 / C28) onto the corresponding new cells, same as using Format
# Painter / Copy > Paste Special > Formats in the Excel UI.
$ws.Range("A" + $prevRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C" + $prevRow).Copy()
$ws.Range("C" + $newRow).PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Column A: running number (one more than the previous row's value, i.e. 28)
$ws.Cells.Item($newRow, 1).Value = 28

# Column C: the new user story text
$ws.Cells.Item($newRow, 3).Value = "AS A Buyer I should be able to define the motivation of my comparision So THAT I can Able to get the result based on my motivation."

# Keep the same row height as the row above it.
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item($prevRow).RowHeight
